$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.849625529173466
$ws.Range("D2").Value = 0.2540673328779519
$ws.Range("E2").Value = 0.1941232997959048
$ws.Range("F2").Value = 4.801729747786169
$ws.Range("G2").Value = 0.002638964791689326
$ws.Range("I2").Value = 1.594364461584895
$ws.Range("J2").Value = 0.214695060966136
$ws.Range("L2").Value = 1.521143295997831
$ws.Range("M2").Value = 0.8262143260414803
$ws.Range("N2").Value = 2.156181415586711
$ws.Range("B3").Value = 1.786001265364916
$ws.Range("D3").Value = 0.2268013807562568
$ws.Range("E3").Value = 0.1693828956727188
$ws.Range("F3").Value = 4.753225963068559
$ws.Range("G3").Value = 0.002648369158398182
$ws.Range("I3").Value = 1.607815780250682
$ws.Range("J3").Value = 0.1865416715813524
$ws.Range("L3").Value = 1.464487373882235
$ws.Range("M3").Value = 0.7966268535794399
$ws.Range("N3").Value = 2.184140424782194
$ws.Range("B4").Value = 1.747910263296717
$ws.Range("D4").Value = 0.2102219068270017
$ws.Range("E4").Value = 0.1542499865882121
$ws.Range("F4").Value = 4.72671600724297
$ws.Range("G4").Value = 0.002654439660905346
$ws.Range("I4").Value = 1.616908952631007
$ws.Range("J4").Value = 0.1692676351267579
$ws.Range("L4").Value = 1.430672333893312
$ws.Range("M4").Value = 0.7789426618601709
$ws.Range("N4").Value = 2.202085626710599
$ws.Range("B5").Value = 1.732631961869487
$ws.Range("D5").Value = 0.203504575911893
$ws.Range("E5").Value = 0.1480964863396252
$ws.Range("F5").Value = 4.716729285004959
$ws.Range("G5").Value = 0.002656988231392118
$ws.Range("I5").Value = 1.620823514763721
$ws.Range("J5").Value = 0.1622302222614138
$ws.Range("L5").Value = 1.417134241713541
$ws.Range("M5").Value = 0.7718566706711272
$ws.Range("N5").Value = 2.209594009119641
$ws.Range("B6").Value = 1.730109726744161
$ws.Range("D6").Value = 0.2023914685844943
$ws.Range("E6").Value = 0.147075468417242
$ws.Range("F6").Value = 4.715120116289654
$ws.Range("G6").Value = 0.002657415946031709
$ws.Range("I6").Value = 1.621486129009476
$ws.Range("J6").Value = 0.1610617455799144
$ws.Range("L6").Value = 1.414900782197122
$ws.Range("M6").Value = 0.7706872969141543
$ws.Range("N6").Value = 2.210852578388955
$ws.Range("B7").Value = 1.747703227704648
$ws.Range("D7").Value = 0.2101311592217314
$ws.Range("E7").Value = 0.1541669461881128
$ws.Range("F7").Value = 4.726578025942416
$ws.Range("G7").Value = 0.002654473728467461
$ws.Range("I7").Value = 1.616960900373748
$ws.Range("J7").Value = 0.1691727199858093
$ws.Range("L7").Value = 1.43048877873315
$ws.Range("M7").Value = 0.778846611135144
$ws.Range("N7").Value = 2.202186095697899
$ws.Range("B8").Value = 1.827484994343365
$ws.Range("D8").Value = 0.2446313367477444
$ws.Range("E8").Value = 0.1855800461987229
$ws.Range("F8").Value = 4.7843227377669
$ws.Range("G8").Value = 0.002642146133931433
$ws.Range("I8").Value = 1.598828954169484
$ws.Range("J8").Value = 0.2049845005963249
$ws.Range("L8").Value = 1.501405015262236
$ws.Range("M8").Value = 0.8159118327697499
$ws.Range("N8").Value = 2.165660076209755
$ws.Range("B9").Value = 1.991725390462818
$ws.Range("D9").Value = 0.3136554460776892
$ws.Range("E9").Value = 0.2476997062652231
$ws.Range("F9").Value = 4.923822175375136
$ws.Range("G9").Value = 0.00262030748124448
$ws.Range("I9").Value = 1.569921900057167
$ws.Range("J9").Value = 0.2753660266536144
$ws.Range("L9").Value = 1.64831126361392
$ws.Range("M9").Value = 0.892472586499423
$ws.Range("N9").Value = 2.100212692780874
$ws.Range("B10").Value = 2.11723631429436
$ws.Range("D10").Value = 0.3653255382764087
$ws.Range("E10").Value = 0.2937436978724293
$ws.Range("F10").Value = 5.042768196509712
$ws.Range("G10").Value = 0.002605666368986848
$ws.Range("I10").Value = 1.552783264960304
$ws.Range("J10").Value = 0.3272551060753983
$ws.Range("L10").Value = 1.761216182363626
$ws.Range("M10").Value = 0.95115743568185
$ws.Range("N10").Value = 2.055901317067277
$ws.Range("B11").Value = 2.175408764755105
$ws.Range("D11").Value = 0.3890671975391342
$ws.Range("E11").Value = 0.3147980224710523
$ws.Range("F11").Value = 5.100553285043986
$ws.Range("G11").Value = 0.002599306211308411
$ws.Range("I11").Value = 1.545886805013019
$ws.Range("J11").Value = 0.3509190559452406
$ws.Range("L11").Value = 1.813703380319794
$ws.Range("M11").Value = 0.9784003121869489
$ws.Range("N11").Value = 2.036563455240295
$ws.Range("B12").Value = 2.197593623203716
$ws.Range("D12").Value = 0.3980937872480013
$ws.Range("E12").Value = 0.322787949298899
$ws.Range("F12").Value = 5.122971690318991
$ws.Range("G12").Value = 0.002596940608545582
$ws.Range("I12").Value = 1.543405584850831
$ws.Range("J12").Value = 0.3598900231349091
$ws.Range("L12").Value = 1.833744306646793
$ws.Range("M12").Value = 0.9887963892160343
$ws.Range("N12").Value = 2.029358761802783
$ws.Range("B13").Value = 2.192808747015647
$ws.Range("D13").Value = 0.3961481083386218
$ws.Range("E13").Value = 0.3210663905852442
$ws.Range("F13").Value = 5.11811949983607
$ws.Range("G13").Value = 0.002597448182702234
$ws.Range("I13").Value = 1.543934149398652
$ws.Range("J13").Value = 0.3579574998046269
$ws.Range("L13").Value = 1.829420736573752
$ws.Range("M13").Value = 0.9865538404607435
$ws.Range("N13").Value = 2.030905162778638
$ws.Range("B14").Value = 2.177230787647147
$ws.Range("D14").Value = 0.3898090840003476
$ws.Range("E14").Value = 0.3154550070035071
$ws.Range("F14").Value = 5.102386863879701
$ws.Range("G14").Value = 0.002599110734739922
$ws.Range("I14").Value = 1.545680056497282
$ws.Range("J14").Value = 0.3516568953406534
$ws.Range("L14").Value = 1.815348830478058
$ws.Range("M14").Value = 0.9792539978110426
$ws.Range("N14").Value = 2.035968351538645
$ws.Range("B15").Value = 2.167709222785504
$ws.Range("D15").Value = 0.3859310181909734
$ws.Range("E15").Value = 0.3120201455331824
$ws.Range("F15").Value = 5.092820267633812
$ws.Range("G15").Value = 0.002600134667211687
$ws.Range("I15").Value = 1.546766473596385
$ws.Range("J15").Value = 0.3477989304918481
$ws.Range("L15").Value = 1.806750988064834
$ws.Range("M15").Value = 0.9747930595213745
$ws.Range("N15").Value = 2.039085093980519
$ws.Range("B16").Value = 2.113456392349008
$ws.Range("D16").Value = 0.3637789063499497
$ws.Range("E16").Value = 0.2923700611208062
$ws.Range("F16").Value = 5.039066440950819
$ws.Range("G16").Value = 0.002606088033719421
$ws.Range("I16").Value = 1.553252151553849
$ws.Range("J16").Value = 0.3257099251236752
$ws.Range("L16").Value = 1.757808955264466
$ws.Range("M16").Value = 0.9493881513917586
$ws.Range("N16").Value = 2.057181607721747
$ws.Range("B17").Value = 2.080450803180725
$ws.Range("D17").Value = 0.350251403844311
$ws.Range("E17").Value = 0.2803442503851983
$ws.Range("F17").Value = 5.007037229389823
$ws.Range("G17").Value = 0.002609816887324529
$ws.Range("I17").Value = 1.557462048991283
$ws.Range("J17").Value = 0.3121752270533023
$ws.Range("L17").Value = 1.728075306547964
$ws.Range("M17").Value = 0.9339439325332961
$ws.Range("N17").Value = 2.068493345928155
$ws.Range("B18").Value = 2.06156813922928
$ws.Range("D18").Value = 0.3424928903565387
$ws.Range("E18").Value = 0.2734374609069903
$ws.Range("F18").Value = 4.988960390136754
$ws.Range("G18").Value = 0.002611989896358051
$ws.Range("I18").Value = 1.559968126844652
$ws.Range("J18").Value = 0.3043959262989233
$ws.Range("L18").Value = 1.711079101645225
$ws.Range("M18").Value = 0.9251122263917892
$ws.Range("N18").Value = 2.075076686303491
$ws.Range("B19").Value = 2.055192147956404
$ws.Range("D19").Value = 0.3398697316345647
$ws.Range("E19").Value = 0.2711006374352252
$ws.Range("F19").Value = 4.982898989703131
$ws.Range("G19").Value = 0.002612730504840385
$ws.Range("I19").Value = 1.56083115437854
$ws.Range("J19").Value = 0.3017628909684618
$ws.Range("L19").Value = 1.70534255399383
$ws.Range("M19").Value = 0.9221307558888867
$ws.Range("N19").Value = 2.077318930761253
$ws.Range("B20").Value = 2.08395381299033
$ws.Range("D20").Value = 0.3516891222353138
$ws.Range("E20").Value = 0.2816233585162848
$ws.Range("F20").Value = 5.010410977826353
$ws.Range("G20").Value = 0.002609417020770782
$ws.Range("I20").Value = 1.557005130168413
$ws.Range("J20").Value = 0.3136154394901496
$ws.Range("L20").Value = 1.731229530949349
$ws.Range("M20").Value = 0.9355826680443329
$ws.Range("N20").Value = 2.067281207699183
$ws.Range("B21").Value = 2.181802161182986
$ws.Range("D21").Value = 0.3916700112983165
$ws.Range("E21").Value = 0.3171027299972025
$ws.Range("F21").Value = 5.106993297853478
$ws.Range("G21").Value = 0.002598621242368575
$ws.Range("I21").Value = 1.545163697125254
$ws.Range("J21").Value = 0.3535072543426168
$ws.Range("L21").Value = 1.819477581769718
$ws.Range("M21").Value = 0.9813959636354781
$ws.Range("N21").Value = 2.034477961833545
$ws.Range("B22").Value = 2.246662924739837
$ws.Range("D22").Value = 0.4180114893977418
$ws.Range("E22").Value = 0.3403910010798938
$ws.Range("F22").Value = 5.173246035432129
$ws.Range("G22").Value = 0.002591815199784051
$ws.Range("I22").Value = 1.538184672570921
$ws.Range("J22").Value = 0.3796374969128351
$ws.Range("L22").Value = 1.878117048433467
$ws.Range("M22").Value = 1.011803248532729
$ws.Range("N22").Value = 2.013727864753888
$ws.Range("B23").Value = 2.211961673041287
$ws.Range("D23").Value = 0.4039324718581554
$ws.Range("E23").Value = 0.3279519217479532
$ws.Range("F23").Value = 5.137596604987607
$ws.Range("G23").Value = 0.002595424975426176
$ws.Range("I23").Value = 1.5418396621114
$ws.Range("J23").Value = 0.3656854543953614
$ws.Range("L23").Value = 1.846730741938018
$ws.Range("M23").Value = 0.9955313091982845
$ws.Range("N23").Value = 2.024739468990596
$ws.Range("B24").Value = 2.082369813651724
$ws.Range("D24").Value = 0.3510390718882945
$ws.Range("E24").Value = 0.2810450520095458
$ws.Range("F24").Value = 5.008884656505558
$ws.Range("G24").Value = 0.002609597709464459
$ws.Range("I24").Value = 1.557211436138779
$ws.Range("J24").Value = 0.3129643133942182
$ws.Range("L24").Value = 1.729803200886977
$ws.Range("M24").Value = 0.9348416482253157
$ws.Range("N24").Value = 2.067828966008776
$ws.Range("B25").Value = 1.946450547531981
$ws.Range("D25").Value = 0.2948234084059322
$ws.Range("E25").Value = 0.2308301352006907
$ws.Range("F25").Value = 4.883226915091768
$ws.Range("G25").Value = 0.002625967443828673
$ws.Range("I25").Value = 1.577025495084328
$ws.Range("J25").Value = 0.2563006783486231
$ws.Range("L25").Value = 1.6077102367706
$ws.Range("M25").Value = 0.8713385335620814
$ws.Range("N25").Value = 2.117255830167051
